# Commit: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# This workbook is an "Estado de Cuenta" (account statement) listing workers
# with overdue periods (2507, 2508). The edit adds a new overdue period
# (2509) for the same three workers, duplicating the existing "2508" block,
# and updates the summary totals (VALOR MORA, Cant. Periodos) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the "2508" data block (rows 19:21) into new rows 22:24 ---
# Insert 3 new rows before row 22 (shifting the signature block, rows 26:27,
# down to rows 29:30), then copy the source block's values AND formatting
# into the newly inserted rows - this is exactly what Excel does when a
# user selects a block, copies it, and inserts the copy below the table.
$src  = $ws.Range("B19:J21")
$dest = $ws.Range("B22:J24")

$dest.Insert(-4121) | Out-Null   # -4121 = xlShiftDown
$src.Copy($dest) | Out-Null

# The source's last row (21) carried the table's special "closing" bottom
# border. After duplication that border belongs on the new last row (24),
# so restore row 21 to the normal interior-row border/format (same as row 20).
$ws.Range("B20:J20").Copy() | Out-Null
$ws.Range("B21:J21").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats

# Set the new block's period ("Periodo Mora") to 2509
$ws.Range("E22").Value = "2509"
$ws.Range("E23").Value = "2509"
$ws.Range("E24").Value = "2509"

# Center-align the "Periodo Mora" column for every data row (16:24) - this
# formatting touch-up was applied across the whole table as part of the edit.
$ws.Range("E16:E24").HorizontalAlignment = -4108   # -4108 = xlCenter

# --- 2. Update summary figures ---
# VALOR MORA (total overdue amount) grows by the new period's total.
$ws.Range("E11").Value = 512460

# Cant. Periodos (count of overdue periods) goes from 2 to 3.
$ws.Range("F13").Value = 3

$excel.CutCopyMode = 0
